$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19
$ws.Range("A19").Value = 111671188
$ws.Range("B19").Value = 78605
$ws.Range("D19").Value = "LC"
$ws.Range("E19").Value = 6462
$ws.Range("F19").Value = "Stuplav"
$ws.Range("G19").Value = "Nephroma bellum"
$ws.Range("H19").Value = "(Spreng.) Tuck."
$ws.Range("Q19").Value = 558215.9329796816
$ws.Range("R19").Value = 7067869.292590594
$ws.Range("L19").ClearContents()

# Row 20
$ws.Range("A20").Value = 111671197
$ws.Range("Q20").Value = 558250.1783714101
$ws.Range("R20").Value = 7067936.828089682

# Row 21
$ws.Range("A21").Value = 111671190
$ws.Range("B21").Value = 78611
$ws.Range("D21").Value = "LC"
$ws.Range("E21").Value = 6463
$ws.Range("F21").Value = "Bårdlav"
$ws.Range("G21").Value = "Nephroma parile"
$ws.Range("H21").Value = "(Ach.) Ach."
$ws.Range("Q21").Value = 558215.9329796816
$ws.Range("R21").Value = 7067869.292590594
$ws.Range("L21").ClearContents()

# Row 22
$ws.Range("A22").Value = 111671201
$ws.Range("Q22").Value = 558250.1783714101
$ws.Range("R22").Value = 7067936.828089682

# Row 23
$ws.Range("A23").Value = 111670567
$ws.Range("B23").Value = 96346
$ws.Range("D23").Value = "NT"
$ws.Range("E23").Value = 620
$ws.Range("F23").Value = "Skogsfru"
$ws.Range("G23").Value = "Epipogium aphyllum"
$ws.Range("H23").Value = "Sw."
$ws.Range("Q23").Value = 558129.9933989302
$ws.Range("R23").Value = 7067958.536170656
$ws.Range("K23").Copy($ws.Range("L23"))

# Row 24
$ws.Range("A24").Value = 111671226
$ws.Range("B24").Value = 78579
$ws.Range("E24").Value = 2081
$ws.Range("F24").Value = "Skrovellav"
$ws.Range("G24").Value = "Lobaria scrobiculata"
$ws.Range("H24").Value = "(Scop.) DC."
$ws.Range("Q24").Value = 558118.4535210516
$ws.Range("R24").Value = 7067742.103054954

# Row 25
$ws.Range("A25").Value = 111670477
$ws.Range("B25").Value = 96346
$ws.Range("D25").Value = "NT"
$ws.Range("E25").Value = 620
$ws.Range("F25").Value = "Skogsfru"
$ws.Range("G25").Value = "Epipogium aphyllum"
$ws.Range("H25").Value = "Sw."
$ws.Range("Q25").Value = 558155.0815836267
$ws.Range("R25").Value = 7068017.481975557
$ws.Range("K25").Copy($ws.Range("L25"))

# Row 26
$ws.Range("A26").Value = 111670558
$ws.Range("Q26").Value = 558133.6011735104
$ws.Range("R26").Value = 7067979.426396712

# Row 27
$ws.Range("A27").Value = 111671294
$ws.Range("B27").Value = 78578
$ws.Range("E27").Value = 6458
$ws.Range("F27").Value = "Lunglav"
$ws.Range("G27").Value = "Lobaria pulmonaria"
$ws.Range("H27").Value = "(L.) Hoffm."
$ws.Range("Q27").Value = 558118.4535210516
$ws.Range("R27").Value = 7067742.103054954
$ws.Range("L27").ClearContents()

# Row 28
$ws.Range("A28").Value = 111670497
$ws.Range("B28").Value = 96346
$ws.Range("E28").Value = 620
$ws.Range("F28").Value = "Skogsfru"
$ws.Range("G28").Value = "Epipogium aphyllum"
$ws.Range("H28").Value = "Sw."
$ws.Range("Q28").Value = 558159.8619213518
$ws.Range("R28").Value = 7068022.886732788
$ws.Range("K28").Copy($ws.Range("L28"))

# Row 29
$ws.Range("A29").Value = 111671179
$ws.Range("B29").Value = 78578
$ws.Range("E29").Value = 6458
$ws.Range("F29").Value = "Lunglav"
$ws.Range("G29").Value = "Lobaria pulmonaria"
$ws.Range("H29").Value = "(L.) Hoffm."
$ws.Range("Q29").Value = 558215.9656782644
$ws.Range("R29").Value = 7067867.520903144
$ws.Range("L29").ClearContents()

# Row 30
$ws.Range("A30").Value = 111670510
$ws.Range("B30").Value = 96346
$ws.Range("E30").Value = 620
$ws.Range("F30").Value = "Skogsfru"
$ws.Range("G30").Value = "Epipogium aphyllum"
$ws.Range("H30").Value = "Sw."
$ws.Range("Q30").Value = 558124.4538526792
$ws.Range("R30").Value = 7067994.321708324
$ws.Range("K30").Copy($ws.Range("L30"))
